# "error solve ifrs list"
# Rewrites the per-year financial figures in the IFRS company_list sheet:
# columns D..AJ for rows 2-6 get new (corrected) figures, the stray V
# column (present only for some of the annual rows) is dropped for
# 2014-2016, zeroed for 2017-2018, and the three forecast rows (7-9,
# 2019E/2020E/2021E) lose all of their figures, leaving only the
# identifying columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12) ---
$ws.Range("D2").Value = 5027
$ws.Range("E2").Value = 541
$ws.Range("F2").Value = 541
$ws.Range("G2").Value = 1636
$ws.Range("H2").Value = 1310
$ws.Range("I2").Value = 1275
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 11758
$ws.Range("L2").Value = 1015
$ws.Range("M2").Value = 10743
$ws.Range("N2").Value = 10559
$ws.Range("O2").Value = 183
$ws.Range("P2").Value = 498
$ws.Range("Q2").Value = 947
$ws.Range("R2").Value = -833
$ws.Range("S2").Value = -561
$ws.Range("T2").Value = 41
$ws.Range("U2").Value = 906
$ws.Range("V2").Value = $null
$ws.Range("W2").Value = 10.77
$ws.Range("X2").Value = 26.06
$ws.Range("Y2").Value = 12.48
$ws.Range("Z2").Value = 11.45
$ws.Range("AA2").Value = 9.449999999999999
$ws.Range("AB2").Value = 1901.79
$ws.Range("AC2").Value = 1279
$ws.Range("AD2").Value = 17.67
$ws.Range("AE2").Value = 10637
$ws.Range("AF2").Value = 2.12
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 2.65
$ws.Range("AI2").Value = 46.71
$ws.Range("AJ2").Value = 99700000

# --- Row 3 (2015/12) ---
$ws.Range("D3").Value = 5094
$ws.Range("E3").Value = 488
$ws.Range("F3").Value = 488
$ws.Range("G3").Value = 1560
$ws.Range("H3").Value = 1250
$ws.Range("I3").Value = 1216
$ws.Range("J3").Value = 34
$ws.Range("K3").Value = 12564
$ws.Range("L3").Value = 1163
$ws.Range("M3").Value = 11401
$ws.Range("N3").Value = 11200
$ws.Range("O3").Value = 201
$ws.Range("P3").Value = 498
$ws.Range("Q3").Value = 1093
$ws.Range("R3").Value = -648
$ws.Range("S3").Value = -613
$ws.Range("T3").Value = 102
$ws.Range("U3").Value = 991
$ws.Range("V3").Value = $null
$ws.Range("W3").Value = 9.59
$ws.Range("X3").Value = 24.53
$ws.Range("Y3").Value = 11.18
$ws.Range("Z3").Value = 10.28
$ws.Range("AA3").Value = 10.2
$ws.Range("AB3").Value = 2029.35
$ws.Range("AC3").Value = 1220
$ws.Range("AD3").Value = 25.54
$ws.Range("AE3").Value = 11282
$ws.Range("AF3").Value = 2.76
$ws.Range("AG3").Value = 670
$ws.Range("AH3").Value = 2.15
$ws.Range("AI3").Value = 54.7
$ws.Range("AJ3").Value = 99700000

# --- Row 4 (2016/12) ---
$ws.Range("D4").Value = 5138
$ws.Range("E4").Value = 456
$ws.Range("F4").Value = 456
$ws.Range("G4").Value = 1522
$ws.Range("H4").Value = 1224
$ws.Range("I4").Value = 1190
$ws.Range("J4").Value = 33
$ws.Range("K4").Value = 13270
$ws.Range("L4").Value = 1296
$ws.Range("M4").Value = 11974
$ws.Range("N4").Value = 11757
$ws.Range("O4").Value = 217
$ws.Range("P4").Value = 498
$ws.Range("Q4").Value = 1065
$ws.Range("R4").Value = -437
$ws.Range("S4").Value = -682
$ws.Range("T4").Value = 43
$ws.Range("U4").Value = 1021
$ws.Range("V4").Value = $null
$ws.Range("W4").Value = 8.880000000000001
$ws.Range("X4").Value = 23.82
$ws.Range("Y4").Value = 10.37
$ws.Range("Z4").Value = 9.470000000000001
$ws.Range("AA4").Value = 10.83
$ws.Range("AB4").Value = 2135.36
$ws.Range("AC4").Value = 1194
$ws.Range("AD4").Value = 23.07
$ws.Range("AE4").Value = 11843
$ws.Range("AF4").Value = 2.33
$ws.Range("AG4").Value = 670
$ws.Range("AH4").Value = 2.43
$ws.Range("AI4").Value = 55.88
$ws.Range("AJ4").Value = 99700000

# --- Row 5 (2017/12) ---
$ws.Range("D5").Value = 5591
$ws.Range("E5").Value = 477
$ws.Range("F5").Value = 477
$ws.Range("G5").Value = 1576
$ws.Range("H5").Value = 1261
$ws.Range("I5").Value = 1223
$ws.Range("J5").Value = 37
$ws.Range("K5").Value = 13877
$ws.Range("L5").Value = 1430
$ws.Range("M5").Value = 12447
$ws.Range("N5").Value = 12202
$ws.Range("O5").Value = 245
$ws.Range("P5").Value = 498
$ws.Range("Q5").Value = 1114
$ws.Range("R5").Value = -277
$ws.Range("S5").Value = -786
$ws.Range("T5").Value = 26
$ws.Range("U5").Value = 1089
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 8.529999999999999
$ws.Range("X5").Value = 22.55
$ws.Range("Y5").Value = 10.21
$ws.Range("Z5").Value = 9.289999999999999
$ws.Range("AA5").Value = 11.49
$ws.Range("AB5").Value = 2254.33
$ws.Range("AC5").Value = 1227
$ws.Range("AD5").Value = 23.71
$ws.Range("AE5").Value = 12329
$ws.Range("AF5").Value = 2.36
$ws.Range("AG5").Value = 700
$ws.Range("AH5").Value = 2.41
$ws.Range("AI5").Value = 56.63
$ws.Range("AJ5").Value = 99700000

# --- Row 6 (2018/12) ---
$ws.Range("D6").Value = 5635
$ws.Range("E6").Value = 432
$ws.Range("F6").Value = 432
$ws.Range("G6").Value = 1538
$ws.Range("H6").Value = 1204
$ws.Range("I6").Value = 1166
$ws.Range("J6").Value = 38
$ws.Range("K6").Value = 14322
$ws.Range("L6").Value = 1420
$ws.Range("M6").Value = 12902
$ws.Range("N6").Value = 12641
$ws.Range("O6").Value = 261
$ws.Range("P6").Value = 498
$ws.Range("Q6").Value = 954
$ws.Range("R6").Value = -265
$ws.Range("S6").Value = -776
$ws.Range("T6").Value = 21
$ws.Range("U6").Value = 933
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 7.67
$ws.Range("X6").Value = 21.37
$ws.Range("Y6").Value = 9.390000000000001
$ws.Range("Z6").Value = 8.539999999999999
$ws.Range("AA6").Value = 11.01
$ws.Range("AB6").Value = 2350.12
$ws.Range("AC6").Value = 1169
$ws.Range("AD6").Value = 15.31
$ws.Range("AE6").Value = 12811
$ws.Range("AF6").Value = 1.4
$ws.Range("AG6").Value = 700
$ws.Range("AH6").Value = 3.91
$ws.Range("AI6").Value = 59.24
$ws.Range("AJ6").Value = 99700000

# --- Rows 7-9 (2019E/2020E/2021E forecast rows): drop every figure,
# keep only the identifying columns A (index), B (annual marker) and C
# (period label).
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
